$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 387, shifting existing rows 387:482 down to 388:483
$ws.Rows.Item(387).Insert()

# Populate the newly inserted row 387 with the new data record
$ws.Cells.Item(387, 1).Value = 3
$ws.Cells.Item(387, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(387, 3).Value = "Coquimbo"
$ws.Cells.Item(387, 4).Value = 44722
$ws.Cells.Item(387, 5).Value = 5
$ws.Cells.Item(387, 6).Value = 100112037
$ws.Cells.Item(387, 7).Value = "Cebollín"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Primera"
$ws.Cells.Item(387, 10).Value = 390
$ws.Cells.Item(387, 11).Value = 6500
$ws.Cells.Item(387, 12).Value = 7000
$ws.Cells.Item(387, 13).Value = 6731
$ws.Cells.Item(387, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(387, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(387, 16).Value = 187
$ws.Cells.Item(387, 17).Value = 36
$ws.Cells.Item(387, 18).Value = "Hortaliza"
